$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data for the Sun Sep  1 20:09:56 UTC 2024 refresh.
# Price values in column D are plain text (e.g. "521.22"), so force the whole
# column to Text format first; otherwise Excel would reinterpret the numeric-
# looking strings as numbers and silently change/reformat them.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 26/27 and 31/32 swap which coin occupies them (Kaspa <-> Binance-PegBSC-USD
# and Fetch.AI <-> Aptos traded places) in addition to their data refreshing.

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "58.353.23"
$ws.Range("E2").Value = "  -0.81%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.487.45"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "521.22"
$ws.Range("E5").Value = "  -2.28%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "135.65"
$ws.Range("E6").Value = "  +0.67%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.33%  "

# Row 8 (XRP)
$ws.Range("D8").Value = "0.560"
$ws.Range("E8").Value = "  -1.40%  "

# Row 9 (LidoStakedEther)
$ws.Range("D9").Value = "2.506.58"
$ws.Range("E9").Value = "  +0.33%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  -1.83%  "

# Row 11 (TRON)
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12 (Toncoin)
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -0.93%  "

# Row 13 (Cardano)
$ws.Range("E13").Value = "  -1.74%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "2.926.93"
$ws.Range("E14").Value = "  -0.28%  "

# Row 15 (WrappedBTC)
$ws.Range("D15").Value = "58.260.63"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16 (Avalanche)
$ws.Range("D16").Value = "22.23"
$ws.Range("E16").Value = "  -1.88%  "

# Row 17 (ShibaInu)
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  -1.46%  "

# Row 18 (WrappedEther)
$ws.Range("D18").Value = "2.489.11"
$ws.Range("E18").Value = "  -0.57%  "

# Row 19 (Chainlink)
$ws.Range("D19").Value = "10.71"
$ws.Range("E19").Value = "  -2.66%  "

# Row 20 (Polkadot)
$ws.Range("D20").Value = "4.20"
$ws.Range("E20").Value = "  -0.85%  "

# Row 21 (BitcoinCash)
$ws.Range("D21").Value = "322.20"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22 (Dai)
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 (Uniswap)
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  -2.59%  "

# Row 24 (Litecoin)
$ws.Range("D24").Value = "64.53"
$ws.Range("E24").Value = "  -0.74%  "

# Row 25 (Polygon)
$ws.Range("D25").Value = "0.415"
$ws.Range("E25").Value = "  -0.93%  "

# Row 26: 'Binance-PegBSC-USD' -> 'Kaspa'
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.161"
$ws.Range("E26").Value = "  -1.03%  "

# Row 27: 'Kaspa' -> 'Binance-PegBSC-USD'
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.57%  "

# Row 28 (InternetComputer(DFINITY))
$ws.Range("D28").Value = "7.42"
$ws.Range("E28").Value = "  -1.02%  "

# Row 29 (PEPE)
$ws.Range("D29").Value = "0.0₃0754"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30 (Monero)
$ws.Range("D30").Value = "169.79"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31: 'Aptos' -> 'Fetch.AI'
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +5.56%  "

# Row 32: 'Fetch.AI' -> 'Aptos'
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.39"
$ws.Range("E32").Value = "  -0.64%  "

# Row 33 (PancakeSwap)
$ws.Range("E33").Value = "  -1.80%  "

# Row 34 (USDe)
$ws.Range("E34").Value = "  -0.07%  "

# Row 35 (FirstDigitalUSD)
$ws.Range("E35").Value = "  -0.38%  "

# Row 36 (EthereumClassic)
$ws.Range("D36").Value = "18.16"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37 (ImmutableX)
$ws.Range("E37").Value = "  -0.69%  "

# Row 38 (NEARProtocol)
$ws.Range("D38").Value = "4.06"
$ws.Range("E38").Value = "  +0.96%  "

# Row 39 (OKB)
$ws.Range("D39").Value = "36.69"
$ws.Range("E39").Value = "  -0.10%  "

# Row 40 (Stacks)
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  -2.33%  "

# Row 41 (SuiNetwork)
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  +1.34%  "

# Row 42 (RenderToken)
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  +5.52%  "

# Row 43 (Bittensor)
$ws.Range("D43").Value = "277.80"
$ws.Range("E43").Value = "  -0.89%  "

# Row 44 (Filecoin)
$ws.Range("D44").Value = "3.47"
$ws.Range("E44").Value = "  -2.52%  "

# Row 45 (Mantle)
$ws.Range("D45").Value = "0.601"
$ws.Range("E45").Value = "  +0.26%  "

# Row 46 (Aave)
$ws.Range("D46").Value = "124.27"
$ws.Range("E46").Value = "  -3.52%  "

# Row 47 (Stellar)
$ws.Range("D47").Value = "0.0912"
$ws.Range("E47").Value = "  -1.41%  "

# Row 48 (Hedera)
$ws.Range("D48").Value = "0.0494"
$ws.Range("E48").Value = "  -0.73%  "

# Row 49 (VeChain)
$ws.Range("D49").Value = "0.0215"
$ws.Range("E49").Value = "  -1.29%  "

# Row 50 (InjectiveProtocol)
$ws.Range("D50").Value = "17.21"
$ws.Range("E50").Value = "  +0.29%  "

# Row 51 (Maker)
$ws.Range("D51").Value = "1.749.36"
$ws.Range("E51").Value = "  +0.17%  "
